$d = $word.ActiveDocument

$pairs = @(
    @("868÷5=", "991÷4="),
    @("694÷9=", "382÷4="),
    @("117÷6=", "173÷3="),
    @("646÷3=", "168÷4="),
    @("726÷3=", "886÷8="),
    @("392÷8=", "458÷2="),
    @("196÷7=", "676÷8="),
    @("941÷3=", "271÷9="),
    @("975÷9=", "788÷7="),
    @("565÷3=", "942÷2="),
    @("585÷6=", "120÷4="),
    @("952÷8=", "586÷6="),
    @("421÷7=", "704÷6="),
    @("534÷6=", "300÷3="),
    @("171÷8=", "192÷2="),
    @("731÷3=", "611÷8="),
    @("773÷8=", "101÷4="),
    @("939÷7=", "522÷4="),
    @("920÷2=", "455÷5="),
    @("375÷9=", "100÷5="),
    @("508÷2=", "175÷9="),
    @("865÷8=", "919÷8="),
    @("887÷4=", "469÷7="),
    @("293÷7=", "343÷2="),
    @("154÷7=", "825÷8=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
